$wb = $excel.ActiveWorkbook

# 1. Rename Sheet2 -> "data pencairan"
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "data pencairan"

# 2. Fill in the data table on the "data pencairan" sheet
# Order of writes matters for shared-string table ordering, so we
# reproduce the same order the original author used: keterangan column
# first, then alamat column, then the header row, then kecamatan/desa.

$kegiatan = @("Lampu Penerangan Jalan ", "Paving Jalan", "Paving Jalan", "Bantuan paving jalan lingkungan", "Pavingisasi Jalan")
$alamat   = @("RT 11 RW 02", "RT 05", "RT 27/ Rw 03 (Jl.Makam)", "Banaran", "RT 08 RW 02")
$tahun    = @(2023, 2023, 2023, 2023, 2023)
$realisasi = @(60000000, 50000000, 50000000, 150000000, 100000000)

# Step 1: column G (keterangan) for rows 2-6
for ($r = 0; $r -lt 5; $r++) {
    $ws.Cells.Item($r + 2, 7).Value = "import excel"
}

# Step 2: column D (alamat) for rows 2-6
for ($r = 0; $r -lt 5; $r++) {
    $ws.Cells.Item($r + 2, 4).Value = $alamat[$r]
}

# Step 3: header row (A1:G1)
$ws.Cells.Item(1, 1).Value = "kecamatan"
$ws.Cells.Item(1, 2).Value = "desa"
$ws.Cells.Item(1, 3).Value = "kegiatan"
$ws.Cells.Item(1, 4).Value = "alamat"
$ws.Cells.Item(1, 5).Value = "tahun"
$ws.Cells.Item(1, 6).Value = "realisasi"
$ws.Cells.Item(1, 7).Value = "keterangan"

# Step 4: columns A/B (kecamatan / desa) for rows 2-6
for ($r = 0; $r -lt 5; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = "Kec A"
    $ws.Cells.Item($r + 2, 2).Value = "Desa A"
}

# Step 5: remaining columns - kegiatan (C), tahun (E), realisasi (F)
for ($r = 0; $r -lt 5; $r++) {
    $ws.Cells.Item($r + 2, 3).Value = $kegiatan[$r]
    $ws.Cells.Item($r + 2, 5).Value = $tahun[$r]
    $ws.Cells.Item($r + 2, 6).Value = $realisasi[$r]
}

Write-Output "done"
